$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on Price cells whose new value would otherwise be
# auto-converted to a number by Excel (e.g. "97.98" -> 97.98 float).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values (Coin, Link, Price, Volume(1h)).
$ws.Range("D2").Value = "44.108.19"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.264.70"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "97.98"
$ws.Range("E5").Value = "  +15.42%  "
$ws.Range("D6").Value = "274.90"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.642"
$ws.Range("E9").Value = "  +8.27%  "
$ws.Range("D10").Value = "47.96"
$ws.Range("E10").Value = "  +7.47%  "
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").Value = "8.23"
$ws.Range("E12").Value = "  +11.50%  "
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "15.70"
$ws.Range("E14").Value = "  +9.37%  "
$ws.Range("D15").Value = "2.594.07"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "0.841"
$ws.Range("E16").Value = "  +7.69%  "
$ws.Range("D17").Value = "2.266.95"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").Value = "44.106.78"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").Value = "70.98"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "10.34"
$ws.Range("E22").Value = "  +15.85%  "
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "235.03"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "11.56"
$ws.Range("E26").Value = "  +8.94%  "
$ws.Range("E27").Value = "  +13.37%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "39.44"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").Value = "3.36"
$ws.Range("E29").Value = "  -6.60%  "
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "173.43"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "0.0924"
$ws.Range("E32").Value = "  +6.76%  "
$ws.Range("D33").Value = "21.29"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  +7.37%  "
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0355"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.42"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "3.62"
$ws.Range("E39").Value = "  +26.48%  "
$ws.Range("D40").Value = "0.254"
$ws.Range("E40").Value = "  +28.50%  "
$ws.Range("D41").Value = "12.55"
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("E42").Value = "  +4.64%  "
$ws.Range("D43").Value = "62.78"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "5.49"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E45").Value = "  +5.10%  "
$ws.Range("D46").Value = "8.56"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("D47").Value = "100.41"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  +5.27%  "
$ws.Range("D49").Value = "1.20"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "0.435"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "2.474.84"
$ws.Range("E51").Value = "  +2.54%  "
